$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 177 (pushes existing rows 177..256 down to 178..257)
$ws.Rows.Item(177).Insert()

# Populate the new row with the new weekly price observation
$ws.Cells.Item(177, 1).Value = 11
$ws.Cells.Item(177, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(177, 3).Value = "Bíobío"
$ws.Cells.Item(177, 4).Value = 45009
$ws.Cells.Item(177, 5).Value = 8
$ws.Cells.Item(177, 6).Value = "Fruta"
$ws.Cells.Item(177, 7).Value = 100108
$ws.Cells.Item(177, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(177, 9).Value = 100108005
$ws.Cells.Item(177, 10).Value = "Piña"
$ws.Cells.Item(177, 11).Value = "Caramelo"
$ws.Cells.Item(177, 12).Value = "Segunda"
$ws.Cells.Item(177, 13).Value = 220
$ws.Cells.Item(177, 14).Value = 23000
$ws.Cells.Item(177, 15).Value = 24000
$ws.Cells.Item(177, 16).Value = 23545
$ws.Cells.Item(177, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(177, 18).Value = "Ecuador"
$ws.Cells.Item(177, 19).Value = 1682
$ws.Cells.Item(177, 20).Value = 14
